$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

$ws.Range("A1").Value = "ulis1111"
$ws.Range("A2").Value = "ulis1122"
$ws.Range("A3").Value = "ulis1133"
$ws.Range("A4").Value = "ulis1144"
$ws.Range("A5").Value = "ulis1154"
$ws.Range("A6").Value = "ulis1164"
$ws.Range("A7").Value = "ulis1174"
$ws.Range("A8").Value = "ulis1184"

$ws.Range("B1").Value = "ubs13"
$ws.Range("B2").Value = "ubs141"
$ws.Range("B3").Value = "ubs152"
$ws.Range("B4").Value = "ubs162"
$ws.Range("B5").Value = "ubs172"
$ws.Range("B6").Value = "ubs182"
$ws.Range("B7").Value = "ubs192"
$ws.Range("B8").Value = "ubs202"

$ws.Range("B8").Select()
